$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: numeric 0, bold + thin border all around + center/top alignment
$r1 = $ws.Range("B1")
$r1.Value = 0
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.VerticalAlignment = -4160
$r1.HorizontalAlignment = -4108

# A2: same formatting as B1 - copy format rather than re-deriving it,
# so both cells land on the identical cellXf entry.
$r2 = $ws.Range("A2")
$r2.Value = 0
$r1.Copy()
$r2.PasteSpecial(-4122)

# B2: plain string label, default style
$ws.Range("B2").Value = "disconnected_elements"
